$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.031.61'
$ws.Range("E2").Value = '  -7.90%  '

$ws.Range("D3").Value = '1.421.61'
$ws.Range("E3").Value = '  -7.68%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '274.15'
$ws.Range("E6").Value = '  -5.53%  '

$ws.Range("E7").Value = '  -3.94%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3078'
$ws.Range("E8").Value = '  -3.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.68'
$ws.Range("E9").Value = '  -7.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.013'
$ws.Range("E10").Value = '  -4.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06616'
$ws.Range("E11").Value = '  -8.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9997'
$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.425'
$ws.Range("E13").Value = '  -3.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.18'
$ws.Range("E14").Value = '  -7.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.170'
$ws.Range("E15").Value = '  -6.52%  '

$ws.Range("D16").Value = '1.421.30'
$ws.Range("E16").Value = '  -7.89%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001009'
$ws.Range("E17").Value = '  -9.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05835'
$ws.Range("E18").Value = '  -11.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.76'
$ws.Range("E19").Value = '  -10.30%  '

$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.655'
$ws.Range("E21").Value = '  -7.97%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.51'
$ws.Range("E22").Value = '  -5.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.02'
$ws.Range("E23").Value = '  +0.42%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.342'
$ws.Range("E24").Value = '  -1.84%  '

$ws.Range("D25").Value = '20.034.52'
$ws.Range("E25").Value = '  -7.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.296'
$ws.Range("E26").Value = '  -3.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '138.91'
$ws.Range("E27").Value = '  -5.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.93'
$ws.Range("E28").Value = '  -7.91%  '

$ws.Range("D29").Value = '1.581.21'
$ws.Range("E29").Value = '  -7.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '109.19'
$ws.Range("E30").Value = '  -7.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.817'
$ws.Range("E31").Value = '  -21.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8915'
$ws.Range("E32").Value = '  -8.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.435'
$ws.Range("E33").Value = '  -8.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07744'
$ws.Range("E34").Value = '  -5.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.473'
$ws.Range("E35").Value = '  -4.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.37'
$ws.Range("E36").Value = '  +6.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.792'
$ws.Range("E37").Value = '  -7.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.000'
$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05687'
$ws.Range("E39").Value = '  -6.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1927'
$ws.Range("E40").Value = '  -5.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02035'
$ws.Range("E41").Value = '  -7.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.087'
$ws.Range("E42").Value = '  -8.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.269'
$ws.Range("E43").Value = '  -14.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5333'
$ws.Range("E44").Value = '  -7.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.538'
$ws.Range("E45").Value = '  -5.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.27'
$ws.Range("E46").Value = '  -6.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5140'
$ws.Range("E47").Value = '  -7.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.806'
$ws.Range("E48").Value = '  -3.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '109.80'
$ws.Range("E49").Value = '  -7.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.049'
$ws.Range("E50").Value = '  -8.32%  '
